$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "16me1033"
$ws.Range("B1").Value = "Emmanuel Menyaga"
$ws.Range("C1").Value = "mathematical sciences"
$ws.Range("D1").Value = 30

$ws.Range("A2").Value = "13ms1023"
$ws.Range("B2").Value = "Ojonugwa Justice Alikali"
$ws.Range("C2").Value = "mathematical sciences"
$ws.Range("D2").Value = 30

$ws.Range("D2").Select()
